$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q3" right before "2022-Q2" and fill it
#    with the fund-holding data for that quarter. We duplicate the
#    "2022-Q2" sheet (so the header / index-column formatting comes
#    along for free) and then overwrite every value.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$refSheet.Copy($refSheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Header row (row 1), columns B..H - keep the text, formatting is already
# carried over from the duplicated sheet.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Data rows (rows 2..10) -> A:index, B:code, C:name, D:size, E:position, F:ratio, G:value, H:rank
# Columns B (fund code) and D..G (size/position/ratio/value) must stay
# textual (e.g. "001959", "4.00", "0.20") - a leading apostrophe forces
# text entry (like typing it in Excel) without leaving behind a custom
# number format, so leading/trailing zeros survive.
$data = @(
    @(0, "'161219", "国投瑞银新兴产业混合（LOF）", "'6.18", "'79.94", "'3.29", "'0.2033", 10),
    @(1, "'161232", "国投瑞银瑞盛灵活配置混合",     "'4.18", "'94.55", "'4.00", "'0.1672", 9),
    @(2, "'001959", "华商乐享互联灵活配置混合A",     "'4.62", "'93.28", "'3.57", "'0.1649", 3),
    @(3, "'000663", "国投瑞银美丽中国灵活配置混合", "'3.85", "'93.40", "'3.63", "'0.1398", 9),
    @(4, "'161225", "国投瑞银瑞盈灵活配置混合（LOF）", "'2.52", "'94.35", "'3.66", "'0.0922", 10),
    @(5, "'013630", "嘉实均衡臻选一年持有期混合A",   "'2.13", "'81.95", "'3.91", "'0.0833", 5),
    @(6, "'013142", "华商乐享互联灵活配置混合C",     "'1.08", "'93.28", "'3.57", "'0.0386", 3),
    @(7, "'013631", "嘉实均衡臻选一年持有期混合C",   "'0.20", "'81.95", "'3.91", "'0.0078", 5),
    @(8, "'001932", "国寿安保灵活优选混合",         "'0.13", "'34.76", "'1.02", "'0.0013", 10)
)

# the duplicated sheet only had 8 data rows (2022-Q2 had 8 funds); this
# quarter has 9, so row 10 needs its index-column style copied in first.
$newSheet.Range("A9").Copy()
$newSheet.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rec = $data[$r]

    $newSheet.Cells.Item($row, 1).Value = $rec[0]
    $newSheet.Cells.Item($row, 2).Value = $rec[1]
    $newSheet.Cells.Item($row, 3).Value = $rec[2]
    $newSheet.Cells.Item($row, 4).Value = $rec[3]
    $newSheet.Cells.Item($row, 5).Value = $rec[4]
    $newSheet.Cells.Item($row, 6).Value = $rec[5]
    $newSheet.Cells.Item($row, 7).Value = $rec[6]
    $newSheet.Cells.Item($row, 8).Value = $rec[7]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: the quarter labels / counts /
#    values shift down by one row and a new 2022-Q3 entry is written in
#    row 2. Row 9 gains the entry that used to be the last row (2020-Q4).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# row 9 is brand new - copy the index-column style down from row 8 first.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$summary = @(
    @(0, "2022-Q3", 9,  0.9),
    @(1, "2022-Q2", 8,  0.76),
    @(2, "2022-Q1", 4,  0.16),
    @(3, "2021-Q4", 17, 4.27),
    @(4, "2021-Q3", 11, 1.28),
    @(5, "2021-Q2", 21, 4.22),
    @(6, "2021-Q1", 26, 2.87),
    @(7, "2020-Q4", 31, 4.01)
)

for ($r = 0; $r -lt $summary.Length; $r++) {
    $row = 2 + $r
    $rec = $summary[$r]
    $total.Cells.Item($row, 1).Value = $rec[0]
    $total.Cells.Item($row, 2).Value = $rec[1]
    $total.Cells.Item($row, 3).Value = $rec[2]
    $total.Cells.Item($row, 4).Value = $rec[3]
}
